# Fill in the "Where Line #'s" columns (G for the left table, N for the
# right table) of the Project 2 check-off sheet with the line numbers /
# flags the student recorded, and move the active selection to Q29.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: this cell is set first so the new shared string "447, 467" lands at
# shared-string index 90, ahead of "not used" (91) and "147 and 148" (92),
# matching the order those strings were authored in.
$ws.Range("G49").Value = "447, 467"

# Left table ("Where Line #'s" = column G), rows 6-50
$ws.Range("G6").Value  = 82
$ws.Range("G7").Value  = 16
$ws.Range("G8").Value  = 75
$ws.Range("G9").Value  = 64
$ws.Range("G10").Value = 65
$ws.Range("G11").Value = 72
$ws.Range("G12").Value = 74
$ws.Range("G13").Value = 64
$ws.Range("G14").Value = 1394
$ws.Range("G16").Value = $True
$ws.Range("G17").Value = $True
$ws.Range("G18").Value = 1406
$ws.Range("G19").Value = $True
$ws.Range("G20").Value = 387
$ws.Range("G21").Value = $True
$ws.Range("G23").Value = 401
$ws.Range("G24").Value = 1796
$ws.Range("G26").Value = $False
$ws.Range("G27").Value = 1796
$ws.Range("G28").Value = $False
$ws.Range("G29").Value = 355
$ws.Range("G30").Value = 174
$ws.Range("G31").Value = 1796
$ws.Range("G32").Value = $False
$ws.Range("G34").Value = 1793
$ws.Range("G35").Value = 2667
$ws.Range("G36").Value = 1764
$ws.Range("G37").Value = 2670
$ws.Range("G38").Value = 2694
$ws.Range("G39").Value = $False
$ws.Range("G40").Value = 2492
$ws.Range("G41").Value = 269
$ws.Range("G42").Value = 363
$ws.Range("G43").Value = 1524
$ws.Range("G45").Value = 2525
$ws.Range("G46").Value = 480
$ws.Range("G47").Value = 630
$ws.Range("G48").Value = 490
$ws.Range("G50").Value = $True

# Right table ("Where Line #'s" = column N), rows 7-29
$ws.Range("N7").Value  = 37
$ws.Range("N8").Value  = 49
$ws.Range("N9").Value  = 545
$ws.Range("N10").Value = 582
$ws.Range("N11").Value = "not used"
$ws.Range("N12").Value = 146
$ws.Range("N13").Value = 39
$ws.Range("N14").Value = 49
$ws.Range("N15").Value = 41
$ws.Range("N16").Value = 460
$ws.Range("N18").Value = 90
$ws.Range("N19").Value = "147 and 148"
$ws.Range("N20").Value = 37
$ws.Range("N21").Value = 388
$ws.Range("N22").Value = 105
$ws.Range("N23").Value = 92
$ws.Range("N24").Value = 109
$ws.Range("N27").Value = 532
$ws.Range("N28").Value = 515
$ws.Range("N29").Value = 477

# Move / record the active selection at Q29 (was Q20).
$ws.Activate()
$ws.Range("Q29").Select()
